# Apply edit: move the "February 11, 2014" _GoBack bookmark's removal,
# and append the older ("4 February 2014" and "28 January 2014") minutes
# blocks after the current minutes, so the most recent entry stays on top
# per the commit message. Finally, re-create the _GoBack bookmark at the
# very end of the document (matching Word's own "last edit" bookmark).

$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark from the "11 February 2014" line.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# 2. Insert the older minutes (4 Feb 2014, then 28 Jan 2014) as raw OOXML
#    right after the last paragraph of the current (11 Feb 2014) minutes,
#    i.e. at the very end of the document body (before the sectPr).
$insertXml = '<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:pict>
          <v:rect id="_x0000_i1025" style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="t" fillcolor="#a0a0a0" stroked="f"/>
        </w:pict>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Date:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> 4 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Feburary</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> 2014</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Attendance:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Rachael Beal</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Sean </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Rappleyea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Minutes:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>We discussed what I had learned from Scott (previous project manager) about having to physically go to Nephi to make any permanent changes. We both thought about the idea of having a local copy of Server Code that we could work from and make official changes periodically throughout the semester (possibly one major update at the end of the semester).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>I also shared the possibility of our team being reassigned to the new UVU Library – Stain Glass Window project. We’re both very excited about that concept.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">We both expressed concerns about keeping Server and Client code separate, since everything we’ve done on </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>PluralSight</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> seems geared towards Client-side programming (which we both enjoy). We then wondered exactly how the other team we’ll be paired up with would fit in with what we have to do. Another class mate was listening in and said that he thinks that the other class would go to the user and get the “feel” and do mock ups of how it would look and feel. Then it would be our responsibility to implement. If we get put on “Roots of Knowledge”, we’d want to make sure that we were in the initial meeting with the other team to make sure they don’t promise anything we can’t really deliver.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>We officially have been assigned “Roots of Knowledge”. We spoke to Luke and think we have a really good concept of what is supposed to happen. I’m going to email Luke for the Client’s information, and Sean is going to email Professor Love for the 3410-Human Factors team information so we can coordinate a group meeting with each other and then the studio.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:pict>
          <v:rect id="_x0000_i1026" style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="t" fillcolor="#a0a0a0" stroked="f"/>
        </w:pict>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Date:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> 28 January 2014</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Attendance:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Rachael Beal</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Sean </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Rappleyea</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Minutes:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>We spent a little bit of time getting to know one another.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Sean is in the Networking track, Rachael in the CS track.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Neither have a lot of prior or additional web development experience outside of the prerequisite Web I class.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>This is going to be a learning process for the both of us.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>We took a look at the Nephi City Website and had a few possible suggestions for improving performance:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Find a better way to load pictures (which currently take way too long)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Keep sites internal or open in a new tab</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>The drop down lists have a very odd look and feel to them</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>We don’t want to make any real decisions until after we’ve talked to the contact on their requirements and desires.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">We asked each other if they had received any email from the teacher with prior group information.  </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Just found out in class that we need to contact Luke. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>We are going to try and have him come to class on Tuesday so we can get the information from him.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>By a flip of a coin, Sean got to pick the Team Leader.  Rachael is the (</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>un)</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>lucky winner.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng = $d.Content
$rng.Collapse(0)
$rng.InsertXML($insertXml)

Write-Output "Done. Paragraph count:"
Write-Output $d.Paragraphs.Count
